# members_data.xlsx cleanup
#
# The sheet contained two rows that shouldn't have been there (an extra
# duplicate entry for "Duraa Zarihun"/"Salamoon Zarihun" and an extra
# duplicate entry for "Walfaanaa Magarsaa"/"Oromiyaa Walfaanaa"), plus a
# "Total money for winners" column (D) that is no longer needed.
#
# This script removes those two rows and the trailing column, shifting the
# remaining data up/left, the same way a user would do it by right-clicking
# the row/column headers in Excel and choosing Delete.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 is "Duraa Zarihun" / 926757080 / 10000 / 120000 - remove it entirely,
# the rows below shift up (row 9 "Salamoon Zarihun" becomes row 8, etc.)
$ws.Rows.Item(8).Delete()

# After the shift above, the old row 12 ("Oromiyaa Walfaanaa" / 912861288 /
# 10000 / 120000) is now row 11 - remove it too.
$ws.Rows.Item(11).Delete()

# Column D ("Total money for winners") is no longer needed - delete the
# whole column, shifting column C's neighbours left.
$ws.Columns.Item(4).Delete()
